$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1) Define the problem</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>a) Do this in your own words.</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>b) What insight can you offer into the problem that is not immediately</w:t></w:r><w:r><w:tab/><w:t>visible</w:t></w:r><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>from</w:t></w:r><w:r><w:tab/></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> word problem </w:t></w:r><w:r><w:t>alone?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>c) What is the</w:t></w:r><w:r><w:tab/><w:t>overall</w:t></w:r><w:r><w:tab/><w:t>goal?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>2) Break the problem apart</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>a) What are the constraints?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>b) What are the sub-goals?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>3) Identify potential solutions</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>a) For each of</w:t></w:r><w:r><w:tab/><w:t>the sub-problems you’ve discussed in #2, what is</w:t></w:r><w:r><w:tab/><w:t>a possible solution?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>4) Evaluate each potential solution</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>a) Does each solution meet the goals?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>b) Will</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">each </w:t></w:r><w:r><w:t>solution wo</w:t></w:r><w:r><w:t xml:space="preserve">rk for ALL </w:t></w:r><w:r><w:t>cases?</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p/><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">5) Choose a solution and develop a plan to implement </w:t></w:r><w:r><w:t>it.</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p/><w:p><w:r><w:t>a</w:t></w:r><w:r><w:t>) Explain the solution in full.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>b) Describe some test cases</w:t></w:r><w:r><w:tab/><w:t>you tried out to make sure it works. (You</w:t></w:r><w:r><w:tab/><w:t>can</w:t></w:r><w:r><w:tab/><w:t>include drawings and diagrams as part of</w:t></w:r><w:r><w:tab/><w:t>your explanation as long as</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">they are clearly communicating the </w:t></w:r><w:r><w:t>solution).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
